$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Build out the flow-chart block G9:J17 with centered text (bulk alignment first,
#    matching a single "center" button click over the whole selection).
$ws.Range("G9:J17").HorizontalAlignment = -4108

# 2) Fill in the chart values top-to-bottom, left-to-right (this is the order the
#    new vocabulary was typed in, which drives shared-string allocation order).
$ws.Range("I9").Value = "Front Gate"
$ws.Range("J9").Value = ""

$ws.Range("H10").Value = "Worker's Entrance"
$ws.Range("I10").Value = "Main Entrance"

$ws.Range("I11").Value = "Intake Desk"
$ws.Range("H11").Value = "|"

$ws.Range("H12").Value = "|"
$ws.Range("I12").Value = "Spindle Room"

$ws.Range("I13").Value = "Cotton Engines"
$ws.Range("J13").Value = "Loading Dock"

$ws.Range("I14").Value = "Central Stairs"
$ws.Range("I14").Interior.Color = 65535

$ws.Range("I15").Value = "Boiler Room"

# 3) Add the standalone "Finishing Room" note further down, with a yellow highlight.
$ws.Range("I22").Value = "Finishing Room"
$ws.Range("I22").Interior.Color = 65535

# 4) Backfill the "Locker Rooms" label in the chart.
$ws.Range("H13").Value = "Locker Rooms"

# 5) Add the "Basement Stairs" header box near the top, also highlighted.
$ws.Range("I4").Value = "Basement Stairs"
$ws.Range("I4").Interior.Color = 65535
$ws.Range("J4").Value = ""

# --- Selection bookkeeping, matches the saved cursor position ---
$ws.Range("K18").Select()
